$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$arr = New-Object 'object[,]' 16,16
$arr[0,0] = 3
$arr[0,1] = 1
$arr[0,2] = 4.987252666666667
$arr[0,3] = 14.961758
$arr[0,4] = 0.07310203833248144
$arr[0,5] = 0.07310203833248144
$arr[0,6] = 3
$arr[0,7] = 1
$arr[0,8] = 3.135398666666667
$arr[0,9] = 9.406196000000001
$arr[0,10] = 0.1723049126704688
$arr[0,11] = 0.1723049126704688
$arr[0,12] = 15.63702536139645
$arr[0,13] = 140.733228252568
$arr[0,14] = 0.01259584033091148
$arr[0,15] = 0.01259584033091148
$arr[1,0] = 3
$arr[1,1] = 1
$arr[1,2] = 4.987252666666667
$arr[1,3] = 14.961758
$arr[1,4] = 0.07310203833248144
$arr[1,5] = 0.07310203833248144
$arr[1,6] = 3
$arr[1,7] = 1
$arr[1,8] = 3.153682
$arr[1,9] = 9.461046
$arr[1,10] = 0.1733096678828815
$arr[1,11] = 0.1733096678828815
$arr[1,12] = 15.72820896431867
$arr[1,13] = 141.553880678868
$arr[1,14] = 0.01266928998496403
$arr[1,15] = 0.01266928998496403
$arr[2,0] = 3
$arr[2,1] = 1
$arr[2,2] = 4.987252666666667
$arr[2,3] = 14.961758
$arr[2,4] = 0.07310203833248144
$arr[2,5] = 0.07310203833248144
$arr[2,6] = 3
$arr[2,7] = 1
$arr[2,8] = 0.4900660000000001
$arr[2,9] = 1.470198
$arr[2,10] = 0.02693143306797965
$arr[2,11] = 0.02693143306797965
$arr[2,12] = 2.444082965342667
$arr[2,13] = 21.996746688084
$arr[2,14] = 0.001968742652484107
$arr[2,15] = 0.001968742652484107
$arr[3,0] = 3
$arr[3,1] = 1
$arr[3,2] = 4.987252666666667
$arr[3,3] = 14.961758
$arr[3,4] = 0.07310203833248144
$arr[3,5] = 0.07310203833248144
$arr[3,6] = 3
$arr[3,7] = 1
$arr[3,8] = 11.417657
$arr[3,9] = 34.252971
$arr[3,10] = 0.62745398637867
$arr[3,11] = 0.6274539863786701
$arr[3,12] = 56.94274032033533
$arr[3,13] = 512.484662883018
$arr[3,14] = 0.04586816536412182
$arr[3,15] = 0.04586816536412183
$arr[4,0] = 3
$arr[4,1] = 1
$arr[4,2] = 54.02320233333333
$arr[4,3] = 162.069607
$arr[4,4] = 0.7918600623966918
$arr[4,5] = 0.7918600623966917
$arr[4,6] = 3
$arr[4,7] = 1
$arr[4,8] = 3.135398666666667
$arr[4,9] = 9.406196000000001
$arr[4,10] = 0.1723049126704688
$arr[4,11] = 0.1723049126704688
$arr[4,12] = 169.3842765649969
$arr[4,13] = 1524.458489084972
$arr[4,14] = 0.136441378898494
$arr[4,15] = 0.1364413788984939
$arr[5,0] = 3
$arr[5,1] = 1
$arr[5,2] = 54.02320233333333
$arr[5,3] = 162.069607
$arr[5,4] = 0.7918600623966918
$arr[5,5] = 0.7918600623966917
$arr[5,6] = 3
$arr[5,7] = 1
$arr[5,8] = 3.153682
$arr[5,9] = 9.461046
$arr[5,10] = 0.1733096678828815
$arr[5,11] = 0.1733096678828815
$arr[5,12] = 170.3720007809913
$arr[5,13] = 1533.348007028922
$arr[5,14] = 0.1372370044236884
$arr[5,15] = 0.1372370044236884
$arr[6,0] = 3
$arr[6,1] = 1
$arr[6,2] = 54.02320233333333
$arr[6,3] = 162.069607
$arr[6,4] = 0.7918600623966918
$arr[6,5] = 0.7918600623966917
$arr[6,6] = 3
$arr[6,7] = 1
$arr[6,8] = 0.4900660000000001
$arr[6,9] = 1.470198
$arr[6,10] = 0.02693143306797965
$arr[6,11] = 0.02693143306797965
$arr[6,12] = 26.47493467468734
$arr[6,13] = 238.274412072186
$arr[6,14] = 0.0213259262696427
$arr[6,15] = 0.02132592626964269
$arr[7,0] = 3
$arr[7,1] = 1
$arr[7,2] = 54.02320233333333
$arr[7,3] = 162.069607
$arr[7,4] = 0.7918600623966918
$arr[7,5] = 0.7918600623966917
$arr[7,6] = 3
$arr[7,7] = 1
$arr[7,8] = 11.417657
$arr[7,9] = 34.252971
$arr[7,10] = 0.62745398637867
$arr[7,11] = 0.6274539863786701
$arr[7,12] = 616.8183942835997
$arr[7,13] = 5551.365548552397
$arr[7,14] = 0.4968557528048667
$arr[7,15] = 0.4968557528048667
$arr[8,0] = 3
$arr[8,1] = 1
$arr[8,2] = 9.123312666666665
$arr[8,3] = 27.369938
$arr[8,4] = 0.1337274842190096
$arr[8,5] = 0.1337274842190096
$arr[8,6] = 3
$arr[8,7] = 1
$arr[8,8] = 3.135398666666667
$arr[8,9] = 9.406196000000001
$arr[8,10] = 0.1723049126704688
$arr[8,11] = 0.1723049126704688
$arr[8,12] = 28.60522237064978
$arr[8,13] = 257.447001335848
$arr[8,14] = 0.02304190248999794
$arr[8,15] = 0.02304190248999794
$arr[9,0] = 3
$arr[9,1] = 1
$arr[9,2] = 9.123312666666665
$arr[9,3] = 27.369938
$arr[9,4] = 0.1337274842190096
$arr[9,5] = 0.1337274842190096
$arr[9,6] = 3
$arr[9,7] = 1
$arr[9,8] = 3.153682
$arr[9,9] = 9.461046
$arr[9,10] = 0.1733096678828815
$arr[9,11] = 0.1733096678828815
$arr[9,12] = 28.77202693723866
$arr[9,13] = 258.948242435148
$arr[9,14] = 0.02317626587680982
$arr[9,15] = 0.02317626587680982
$arr[10,0] = 3
$arr[10,1] = 1
$arr[10,2] = 9.123312666666665
$arr[10,3] = 27.369938
$arr[10,4] = 0.1337274842190096
$arr[10,5] = 0.1337274842190096
$arr[10,6] = 3
$arr[10,7] = 1
$arr[10,8] = 0.4900660000000001
$arr[10,9] = 1.470198
$arr[10,10] = 0.02693143306797965
$arr[10,11] = 0.02693143306797965
$arr[10,12] = 4.471025345302666
$arr[10,13] = 40.239228107724
$arr[10,14] = 0.003601472790593561
$arr[10,15] = 0.003601472790593561
$arr[11,0] = 3
$arr[11,1] = 1
$arr[11,2] = 9.123312666666665
$arr[11,3] = 27.369938
$arr[11,4] = 0.1337274842190096
$arr[11,5] = 0.1337274842190096
$arr[11,6] = 3
$arr[11,7] = 1
$arr[11,8] = 11.417657
$arr[11,9] = 34.252971
$arr[11,10] = 0.62745398637867
$arr[11,11] = 0.6274539863786701
$arr[11,12] = 104.1668547317553
$arr[11,13] = 937.501692585798
$arr[11,14] = 0.08390784306160824
$arr[11,15] = 0.08390784306160826
$arr[12,0] = 3
$arr[12,1] = 1
$arr[12,2] = 0.08940066666666667
$arr[12,3] = 0.268202
$arr[12,4] = 0.001310415051817319
$arr[12,5] = 0.001310415051817319
$arr[12,6] = 3
$arr[12,7] = 1
$arr[12,8] = 3.135398666666667
$arr[12,9] = 9.406196000000001
$arr[12,10] = 0.1723049126704688
$arr[12,11] = 0.1723049126704688
$arr[12,12] = 0.2803067310657778
$arr[12,13] = 2.522760579592001
$arr[12,14] = 0.0002257909510654511
$arr[12,15] = 0.000225790951065451
$arr[13,0] = 3
$arr[13,1] = 1
$arr[13,2] = 0.08940066666666667
$arr[13,3] = 0.268202
$arr[13,4] = 0.001310415051817319
$arr[13,5] = 0.001310415051817319
$arr[13,6] = 3
$arr[13,7] = 1
$arr[13,8] = 3.153682
$arr[13,9] = 9.461046
$arr[13,10] = 0.1733096678828815
$arr[13,11] = 0.1733096678828815
$arr[13,12] = 0.2819412732546667
$arr[13,13] = 2.537471459292
$arr[13,14] = 0.0002271075974191885
$arr[13,15] = 0.0002271075974191885
$arr[14,0] = 3
$arr[14,1] = 1
$arr[14,2] = 0.08940066666666667
$arr[14,3] = 0.268202
$arr[14,4] = 0.001310415051817319
$arr[14,5] = 0.001310415051817319
$arr[14,6] = 3
$arr[14,7] = 1
$arr[14,8] = 0.4900660000000001
$arr[14,9] = 1.470198
$arr[14,10] = 0.02693143306797965
$arr[14,11] = 0.02693143306797965
$arr[14,12] = 0.04381222711066667
$arr[14,13] = 0.394310043996
$arr[14,14] = 0.00003529135525929122
$arr[14,15] = 0.00003529135525929121
$arr[15,0] = 3
$arr[15,1] = 1
$arr[15,2] = 0.08940066666666667
$arr[15,3] = 0.268202
$arr[15,4] = 0.001310415051817319
$arr[15,5] = 0.001310415051817319
$arr[15,6] = 3
$arr[15,7] = 1
$arr[15,8] = 11.417657
$arr[15,9] = 34.252971
$arr[15,10] = 0.62745398637867
$arr[15,11] = 0.6274539863786701
$arr[15,12] = 1.020746147571333
$arr[15,13] = 9.186715328142
$arr[15,14] = 0.0008222251480733884
$arr[15,15] = 0.0008222251480733884
$ws.Range("E2:T17").Value = $arr
